$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 664, shifting the old
# rows 664:669 down to 668:673 (formatting is carried along automatically).
$ws.Rows("664:667").Insert()

# Row 664 - new entry: "Cuatro cascos verde" Primera
$ws.Cells.Item(664,1).Value  = 4
$ws.Cells.Item(664,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(664,3).Value  = "Los Lagos"
$ws.Cells.Item(664,4).Value  = "2022-04-05"
$ws.Cells.Item(664,5).Value  = 10
$ws.Cells.Item(664,6).Value  = 100112002
$ws.Cells.Item(664,7).Value  = "Pimiento"
$ws.Cells.Item(664,8).Value  = "Cuatro cascos verde"
$ws.Cells.Item(664,9).Value  = "Primera"
$ws.Cells.Item(664,10).Value = 180
$ws.Cells.Item(664,11).Value = 17000
$ws.Cells.Item(664,12).Value = 17000
$ws.Cells.Item(664,13).Value = 17000
$ws.Cells.Item(664,14).Value = "$/caja 18 kilos"
$ws.Cells.Item(664,15).Value = "Provincia de Limarí"
$ws.Cells.Item(664,16).Value = 944
$ws.Cells.Item(664,17).Value = 18
$ws.Cells.Item(664,18).Value = "Hortaliza"

# Row 665 - new entry: "Morrón rojo" Primera
$ws.Cells.Item(665,1).Value  = 4
$ws.Cells.Item(665,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(665,3).Value  = "Los Lagos"
$ws.Cells.Item(665,4).Value  = "2022-04-05"
$ws.Cells.Item(665,5).Value  = 10
$ws.Cells.Item(665,6).Value  = 100112002
$ws.Cells.Item(665,7).Value  = "Pimiento"
$ws.Cells.Item(665,8).Value  = "Morrón rojo"
$ws.Cells.Item(665,9).Value  = "Primera"
$ws.Cells.Item(665,10).Value = 80
$ws.Cells.Item(665,11).Value = 20000
$ws.Cells.Item(665,12).Value = 20000
$ws.Cells.Item(665,13).Value = 20000
$ws.Cells.Item(665,14).Value = "$/caja 20 kilos"
$ws.Cells.Item(665,15).Value = "Provincia de Limarí"
$ws.Cells.Item(665,16).Value = 1000
$ws.Cells.Item(665,17).Value = 20
$ws.Cells.Item(665,18).Value = "Hortaliza"

# Row 666 - new entry: "Morrón rojo" Segunda
$ws.Cells.Item(666,1).Value  = 4
$ws.Cells.Item(666,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(666,3).Value  = "Los Lagos"
$ws.Cells.Item(666,4).Value  = "2022-04-05"
$ws.Cells.Item(666,5).Value  = 10
$ws.Cells.Item(666,6).Value  = 100112002
$ws.Cells.Item(666,7).Value  = "Pimiento"
$ws.Cells.Item(666,8).Value  = "Morrón rojo"
$ws.Cells.Item(666,9).Value  = "Segunda"
$ws.Cells.Item(666,10).Value = 80
$ws.Cells.Item(666,11).Value = 18000
$ws.Cells.Item(666,12).Value = 18000
$ws.Cells.Item(666,13).Value = 18000
$ws.Cells.Item(666,14).Value = "$/caja 20 kilos"
$ws.Cells.Item(666,15).Value = "Provincia de Limarí"
$ws.Cells.Item(666,16).Value = 900
$ws.Cells.Item(666,17).Value = 20
$ws.Cells.Item(666,18).Value = "Hortaliza"

# Row 667 - new entry: "Morrón rojo" Tercera
$ws.Cells.Item(667,1).Value  = 4
$ws.Cells.Item(667,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(667,3).Value  = "Los Lagos"
$ws.Cells.Item(667,4).Value  = "2022-04-05"
$ws.Cells.Item(667,5).Value  = 10
$ws.Cells.Item(667,6).Value  = 100112002
$ws.Cells.Item(667,7).Value  = "Pimiento"
$ws.Cells.Item(667,8).Value  = "Morrón rojo"
$ws.Cells.Item(667,9).Value  = "Tercera"
$ws.Cells.Item(667,10).Value = 80
$ws.Cells.Item(667,11).Value = 15000
$ws.Cells.Item(667,12).Value = 15000
$ws.Cells.Item(667,13).Value = 15000
$ws.Cells.Item(667,14).Value = "$/caja 20 kilos"
$ws.Cells.Item(667,15).Value = "Provincia de Limarí"
$ws.Cells.Item(667,16).Value = 750
$ws.Cells.Item(667,17).Value = 20
$ws.Cells.Item(667,18).Value = "Hortaliza"
